$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'322.03"
$ws.Range("E2").Value = "'-2.85%"
$ws.Range("D3").Value = "'42.59"
$ws.Range("E3").Value = "'-6.25%"
$ws.Range("D4").Value = "'5.159"
$ws.Range("E4").Value = "'-8.07%"
$ws.Range("D5").Value = "'0.08196"
$ws.Range("E5").Value = "'-1.85%"
$ws.Range("D6").Value = "'4.293"
$ws.Range("E6").Value = "'-3.36%"
$ws.Range("D7").Value = "'1.802"
$ws.Range("E7").Value = "'-13.28%"
$ws.Range("E8").Value = "'-3.57%"
$ws.Range("D9").Value = "'0.1108"
$ws.Range("E9").Value = "'-5.56%"
$ws.Range("D10").Value = "'0.1867"
$ws.Range("E10").Value = "'-2.98%"
$ws.Range("D11").Value = "'0.09464"
$ws.Range("E11").Value = "'-3.86%"
$ws.Range("D12").Value = "'0.04639"
$ws.Range("E12").Value = "'0.45%"
$ws.Range("D13").Value = "'7.409"
$ws.Range("E13").Value = "'-28.06%"
$ws.Range("D14").Value = "'0.1058"
$ws.Range("E14").Value = "'-0.35%"
$ws.Range("D15").Value = "'0.001291"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("D16").Value = "'0.005844"
$ws.Range("E16").Value = "'-4.49%"
$ws.Range("E17").Value = "'-0.36%"
$ws.Range("E18").Value = "'-0.98%"
$ws.Range("D19").Value = "'0.3375"
$ws.Range("E19").Value = "'0.79%"
$ws.Range("E20").Value = "'-0.22%"
$ws.Range("E21").Value = "'-12.49%"
$ws.Range("D22").Value = "'0.04164"
$ws.Range("E22").Value = "'-0.37%"
$ws.Range("E23").Value = "'-5.34%"
$ws.Range("D24").Value = "'0.004436"
$ws.Range("E24").Value = "'-2.67%"
$ws.Range("E25").Value = "'-7.95%"
$ws.Range("D26").Value = "'0.0002980"
$ws.Range("E26").Value = "'-20.57%"
$ws.Range("D38").Value = "'0.02772"
$ws.Range("E38").Value = "'2.39%"
$ws.Range("D39").Value = "'0.05595"
$ws.Range("E39").Value = "'-2.84%"
$ws.Range("D40").Value = "'0.008065"
$ws.Range("E40").Value = "'2.10%"
$ws.Range("D41").Value = "'0.1397"
$ws.Range("E41").Value = "'-2.53%"
$ws.Range("D42").Value = "'0.006552"
$ws.Range("E42").Value = "'-9.76%"
$ws.Range("D43").Value = "'0.002086"
$ws.Range("E43").Value = "'2.94%"
$ws.Range("D44").Value = "'0.008332"
$ws.Range("E44").Value = "'-8.59%"
$ws.Range("D45").Value = "'0.3490"
$ws.Range("E45").Value = "'-1.54%"
$ws.Range("D46").Value = "'0.00006957"
$ws.Range("E46").Value = "'-2.49%"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("D48").Value = "'0.003476"
$ws.Range("E48").Value = "'-0.54%"
$ws.Range("E49").Value = "'0.63%"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("E51").Value = "'-0.26%"
